$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F10").Value = 4983
$ws.Range("F14").Value = 1397
$ws.Range("F16").Value = 557
$ws.Range("F17").Value = 6806
$ws.Range("F21").Value = 4585
$ws.Range("F25").Value = 2231
$ws.Range("F29").Value = 182
$ws.Range("F32").Value = 133
$ws.Range("F34").Value = 1250
$ws.Range("F35").Value = 1959
$ws.Range("F36").Value = 200
$ws.Range("F39").Value = 1327
$ws.Range("F44").Value = 1064
$ws.Range("F45").Value = 1347
$ws.Range("F48").Value = 219

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 243
$ws.Range("F20").Value = 126

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1628
$ws.Range("F7").Value = 526
$ws.Range("F10").Value = 1701
$ws.Range("F11").Value = 2050
$ws.Range("F12").Value = 516

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1628
$ws.Range("F4").Value = 526
$ws.Range("F8").Value = 2050
$ws.Range("F9").Value = 4983
$ws.Range("F10").Value = 516
$ws.Range("F13").Value = 892
$ws.Range("F16").Value = 1397
$ws.Range("F18").Value = 557
$ws.Range("F19").Value = 6806
$ws.Range("F24").Value = 4585
$ws.Range("F27").Value = 2231
$ws.Range("F31").Value = 182
$ws.Range("F33").Value = 243
$ws.Range("F35").Value = 133
$ws.Range("F37").Value = 1959
$ws.Range("F38").Value = 200
$ws.Range("F42").Value = 1327
$ws.Range("F47").Value = 1064
$ws.Range("F48").Value = 1347
$ws.Range("F49").Value = 219
